$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so the date-formatted-looking strings are stored
# as literal text (matching the original inline-string/text cell type)
# instead of being auto-converted into Excel date serial numbers.
$ws.Range("D2").Value = "'2026-02-12"
$ws.Range("D3").Value = "'2026-02-13"
$ws.Range("D4").Value = "'2026-02-14"
$ws.Range("D5").Value = "'2026-02-15"
$ws.Range("D6").Value = "'2026-02-16"
